# Add new columns I (I0) and J (IF) to the worksheet, matching the style
# of the existing header row and filling in the per-row numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add I1 = "I0" and J1 = "IF" ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the style from the existing header cell (H1) so the new headers
# match the bold/bordered/centered look of the other header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Data rows 2-26: add values for columns I and J ---
$data = @(
    @{Row=2;  I=1;  J=5}
    @{Row=3;  I=1;  J=5}
    @{Row=4;  I=1;  J=4}
    @{Row=5;  I=1;  J=6}
    @{Row=6;  I=1;  J=5}
    @{Row=7;  I=1;  J=6}
    @{Row=8;  I=1;  J=6}
    @{Row=9;  I=1;  J=6}
    @{Row=10; I=1;  J=7}
    @{Row=11; I=1;  J=5}
    @{Row=12; I=3;  J=8}
    @{Row=13; I=1;  J=4}
    @{Row=14; I=1;  J=5}
    @{Row=15; I=8;  J=9}
    @{Row=16; I=8;  J=9}
    @{Row=17; I=6;  J=6}
    @{Row=18; I=6;  J=8}
    @{Row=19; I=7;  J=8}
    @{Row=20; I=7;  J=8}
    @{Row=21; I=5;  J=7}
    @{Row=22; I=11; J=11}
    @{Row=23; I=7;  J=7}
    @{Row=24; I=7;  J=8}
    @{Row=25; I=8;  J=8}
    @{Row=26; I=4;  J=4}
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 9).Value  = $entry.I   # column I
    $ws.Cells.Item($r, 10).Value = $entry.J   # column J
}
